# Actualización desde MV -datos-
# Append five new daily rows (20-09-2021 .. 24-09-2021) to the "Diaria" sheet,
# following the existing layout: A=Serie (date), C=1 año, D=2 años, E=5 años.
# (Column B is not used on this sheet, same as the pre-existing rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Serie = "20-09-2021"; C = 4.66; D = 3.57; E = 3.23 },
    @{ Serie = "21-09-2021"; C = 4.5;           E = 3.26 },
    @{ Serie = "22-09-2021";                    E = 3.24 },
    @{ Serie = "23-09-2021"; C = 4.5;  D = 3.37; E = 3.32 },
    @{ Serie = "24-09-2021"; C = 4.52; D = 3.42; E = 3.27 }
)

$firstNewRow = 180

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $firstNewRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.Serie

    if ($data.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $data.C }
    if ($data.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $data.D }
    if ($data.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $data.E }
}
